$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.923.41'
$ws.Range('E2').Value = '  +5.82%  '
$ws.Range('D3').Value = '2.233.85'
$ws.Range('E3').Value = '  +3.05%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('E5').Value = '  +2.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.622'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.81'
$ws.Range('E7').Value = '  -2.00%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  +2.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.31'
$ws.Range('E10').Value = '  +1.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0890'
$ws.Range('E11').Value = '  +4.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.103'
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').Value = '2.566.61'
$ws.Range('E13').Value = '  +3.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.65'
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.00'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.802'
$ws.Range('E16').Value = '  -0.93%  '
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').Value = '2.253.38'
$ws.Range('E18').Value = '  +3.88%  '
$ws.Range('D19').Value = '41.753.29'
$ws.Range('E19').Value = '  +5.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.21'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').Value = '0.0₃0897'
$ws.Range('E21').Value = '  -1.53%  '
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.88'
$ws.Range('E23').Value = '  +9.74%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.39'
$ws.Range('E25').Value = '  +2.15%  '
$ws.Range('E26').Value = '  +1.83%  '
$ws.Range('E27').Value = '  +2.17%  '
$ws.Range('E28').Value = '  +2.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.68'
$ws.Range('E29').Value = '  -2.42%  '
$ws.Range('E30').Value = '  +1.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.41'
$ws.Range('E31').Value = '  -1.69%  '
$ws.Range('E32').Value = '  -1.55%  '
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.00'
$ws.Range('E34').Value = '  +6.15%  '
$ws.Range('E35').Value = '  +3.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0635'
$ws.Range('E36').Value = '  +3.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.65'
$ws.Range('E37').Value = '  -4.54%  '
$ws.Range('E38').Value = '  -3.86%  '
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.000254'
$ws.Range('E40').Value = '  +30.29%  '
$ws.Range('E42').Value = '  +5.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.84'
$ws.Range('E43').Value = '  -2.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.58'
$ws.Range('E44').Value = '  +8.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0980'
$ws.Range('E45').Value = '  +6.33%  '
$ws.Range('E46').Value = '  +0.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '98.87'
$ws.Range('E47').Value = '  -3.55%  '
$ws.Range('D48').Value = '1.479.88'
$ws.Range('E48').Value = '  -2.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.49'
$ws.Range('E49').Value = '  -6.62%  '
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.52'
$ws.Range('E51').Value = '  +8.40%  '
